$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing F/G values for rows 624-719 ---
$ws.Cells.Item(624, 6).Value = 51677
$ws.Cells.Item(624, 7).Value = 3965
$ws.Cells.Item(627, 6).Value = 34134
$ws.Cells.Item(627, 7).Value = 2763
$ws.Cells.Item(630, 6).Value = 46889
$ws.Cells.Item(630, 7).Value = 2976
$ws.Cells.Item(638, 6).Value = 37767
$ws.Cells.Item(641, 6).Value = 34417
$ws.Cells.Item(641, 7).Value = 1390
$ws.Cells.Item(642, 6).Value = 67512
$ws.Cells.Item(642, 7).Value = 2390
$ws.Cells.Item(650, 6).Value = 38069
$ws.Cells.Item(671, 6).Value = 32634
$ws.Cells.Item(677, 6).Value = 56199
$ws.Cells.Item(679, 6).Value = 29481
$ws.Cells.Item(681, 6).Value = 26436
$ws.Cells.Item(681, 7).Value = 580
$ws.Cells.Item(684, 6).Value = 57270
$ws.Cells.Item(685, 6).Value = 34487
$ws.Cells.Item(686, 6).Value = 34448
$ws.Cells.Item(687, 6).Value = 31492
$ws.Cells.Item(688, 6).Value = 32162
$ws.Cells.Item(688, 7).Value = 1351
$ws.Cells.Item(689, 6).Value = 15789
$ws.Cells.Item(689, 7).Value = 1066
$ws.Cells.Item(690, 6).Value = 27830
$ws.Cells.Item(690, 7).Value = 1546
$ws.Cells.Item(691, 6).Value = 62611
$ws.Cells.Item(691, 7).Value = 2820
$ws.Cells.Item(692, 6).Value = 41670
$ws.Cells.Item(692, 7).Value = 2687
$ws.Cells.Item(693, 6).Value = 39589
$ws.Cells.Item(693, 7).Value = 2730
$ws.Cells.Item(694, 6).Value = 37611
$ws.Cells.Item(694, 7).Value = 2777
$ws.Cells.Item(695, 6).Value = 37278
$ws.Cells.Item(695, 7).Value = 3134
$ws.Cells.Item(696, 6).Value = 17829
$ws.Cells.Item(696, 7).Value = 2218
$ws.Cells.Item(697, 6).Value = 28926
$ws.Cells.Item(697, 7).Value = 3037
$ws.Cells.Item(698, 6).Value = 70786
$ws.Cells.Item(698, 7).Value = 5821
$ws.Cells.Item(699, 6).Value = 43541
$ws.Cells.Item(699, 7).Value = 4304
$ws.Cells.Item(700, 6).Value = 43713
$ws.Cells.Item(700, 7).Value = 4317
$ws.Cells.Item(701, 6).Value = 41835
$ws.Cells.Item(701, 7).Value = 3851
$ws.Cells.Item(702, 6).Value = 36405
$ws.Cells.Item(702, 7).Value = 3935
$ws.Cells.Item(703, 6).Value = 17104
$ws.Cells.Item(703, 7).Value = 2611
$ws.Cells.Item(704, 6).Value = 25027
$ws.Cells.Item(704, 7).Value = 3703
$ws.Cells.Item(705, 6).Value = 55963
$ws.Cells.Item(705, 7).Value = 6304
$ws.Cells.Item(706, 6).Value = 40651
$ws.Cells.Item(706, 7).Value = 4952
$ws.Cells.Item(707, 6).Value = 38787
$ws.Cells.Item(707, 7).Value = 4615
$ws.Cells.Item(708, 6).Value = 35510
$ws.Cells.Item(708, 7).Value = 4145
$ws.Cells.Item(709, 6).Value = 32329
$ws.Cells.Item(709, 7).Value = 3971
$ws.Cells.Item(710, 6).Value = 14669
$ws.Cells.Item(710, 7).Value = 2629
$ws.Cells.Item(711, 6).Value = 22473
$ws.Cells.Item(711, 7).Value = 3801
$ws.Cells.Item(712, 6).Value = 51230
$ws.Cells.Item(712, 7).Value = 6293
$ws.Cells.Item(713, 6).Value = 37062
$ws.Cells.Item(713, 7).Value = 4733
$ws.Cells.Item(714, 6).Value = 32297
$ws.Cells.Item(714, 7).Value = 3965
$ws.Cells.Item(715, 6).Value = 31639
$ws.Cells.Item(715, 7).Value = 3547
$ws.Cells.Item(716, 6).Value = 29537
$ws.Cells.Item(716, 7).Value = 3656
$ws.Cells.Item(717, 6).Value = 12415
$ws.Cells.Item(717, 7).Value = 2105
$ws.Cells.Item(718, 6).Value = 16759
$ws.Cells.Item(718, 7).Value = 2781
$ws.Cells.Item(719, 6).Value = 42792
$ws.Cells.Item(719, 7).Value = 5116

# --- Add F/G values for rows 720-722 (previously blank) ---
$ws.Cells.Item(720, 6).Value = 30479
$ws.Cells.Item(720, 7).Value = 3433
$ws.Cells.Item(721, 6).Value = 27008
$ws.Cells.Item(721, 7).Value = 3007
$ws.Cells.Item(722, 6).Value = 26039
$ws.Cells.Item(722, 7).Value = 2711

# --- Add new rows 723-725 (new daily stats through 2022-02-27) ---
$ws.Cells.Item(723, 1).Value = 44617
$ws.Cells.Item(723, 2).Value = 1433930
$ws.Cells.Item(723, 3).Value = 23621
$ws.Cells.Item(723, 4).Value = 12384
$ws.Cells.Item(723, 5).Value = 18413
$ws.Cells.Item(723, 6).Value = 18636
$ws.Cells.Item(723, 7).Value = 2356
$ws.Cells.Item(724, 1).Value = 44618
$ws.Cells.Item(724, 2).Value = 1442194
$ws.Cells.Item(724, 3).Value = 15299
$ws.Cells.Item(724, 4).Value = 8264
$ws.Cells.Item(724, 5).Value = 18453
$ws.Cells.Item(724, 6).Value = 7486
$ws.Cells.Item(724, 7).Value = 1250
$ws.Cells.Item(725, 1).Value = 44619
$ws.Cells.Item(725, 2).Value = 1446922
$ws.Cells.Item(725, 3).Value = 8825
$ws.Cells.Item(725, 4).Value = 4728
$ws.Cells.Item(725, 5).Value = 18485
$ws.Cells.Item(725, 6).Value = 7571
$ws.Cells.Item(725, 7).Value = 1421
